$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "52.185.01"
$ws.Range("E2").Value = "  +0.98%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.000.04"
$ws.Range("E3").Value = "  +2.12%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "354.19"
$ws.Range("E5").Value = "  +0.34%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "107.36"
$ws.Range("E6").Value = "  -3.19%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.558"
$ws.Range("E7").Value = "  -0.64%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.613"
$ws.Range("E9").Value = "  -1.84%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.24"
$ws.Range("E10").Value = "  -2.33%  "
$ws.Range("E11").Value = "  +2.14%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0856"
$ws.Range("E12").Value = "  -4.74%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.11"
$ws.Range("E13").Value = "  -3.29%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.452.92"
$ws.Range("E14").Value = "  +1.39%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.59"
$ws.Range("E15").Value = "  -5.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.006.83"
$ws.Range("E16").Value = "  +2.61%  "
$ws.Range("E17").Value = "  +3.82%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "52.213.58"
$ws.Range("E18").Value = "  +0.84%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.42"
$ws.Range("E19").Value = "  +5.95%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.52"
$ws.Range("E20").Value = "  -1.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.59"
$ws.Range("E21").Value = "  -5.49%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0975"
$ws.Range("E22").Value = "  -1.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.19"
$ws.Range("E23").Value = "  -2.65%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "264.47"
$ws.Range("E24").Value = "  -2.73%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.73"
$ws.Range("E25").Value = "  -1.59%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.177"
$ws.Range("E26").Value = "  -1.56%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "26.86"
$ws.Range("E27").Value = "  -1.25%  "
$ws.Range("B28").Value = "Filecoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.54"
$ws.Range("E28").Value = "  -2.62%  "
$ws.Range("E29").Value = "  -0.20%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.105"
$ws.Range("E30").Value = "  -2.48%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.53"
$ws.Range("E31").Value = "  +4.53%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.28"
$ws.Range("E32").Value = "  -3.39%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "36.11"
$ws.Range("E33").Value = "  -2.44%  "
$ws.Range("E34").Value = "  +15.84%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "51.01"
$ws.Range("E35").Value = "  -2.89%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0437"
$ws.Range("E36").Value = "  -0.96%  "
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.28"
$ws.Range("E38").Value = "  -1.14%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.99"
$ws.Range("E39").Value = "  -2.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.73"
$ws.Range("E40").Value = "  -4.35%  "
$ws.Range("E41").Value = "  +1.70%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.117"
$ws.Range("E42").Value = "  -0.59%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "22.84"
$ws.Range("E43").Value = "  -2.21%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "123.55"
$ws.Range("E44").Value = "  +8.58%  "
$ws.Range("E45").Value = "  -0.73%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.122.76"
$ws.Range("E46").Value = "  -1.27%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.35"
$ws.Range("E47").Value = "  -4.04%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.37"
$ws.Range("E48").Value = "  -6.28%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.297.35"
$ws.Range("E49").Value = "  +2.16%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.244"
$ws.Range("E50").Value = "  -0.27%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0333"
$ws.Range("E51").Value = "  +0.99%  "
